# Update the "Google_Search_Test_2" worksheet: cell A2 value changes
# from "LinkedIn" to "Cigniti".
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Google_Search_Test_2")
$ws.Range("A2").Value = "Cigniti"
